# Adds the new "DIETARY_ASSESS_INSTR" variable (Variables sheet) and its
# category/response options (Categories sheet), per the commit:
# "adding information on the variable DIETARY_ASSESS_INSTR for Variables
#  and Categories sheet in both P1 and P2 Dataschemas"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Variables")
$ws2 = $wb.Worksheets.Item("Categories")

# --- Variables sheet: append the new variable as row 129 ---
$ws1.Range("A129").Value = 128
$ws1.Range("B129").Value = "DIETARY_ASSESS_INSTR"
$ws1.Range("C129").Value = "Dietary Assessment Instrument"
$ws1.Range("D129").Value = "integer"

# --- Categories sheet: append the category codes for DIETARY_ASSESS_INSTR ---
$ws2.Range("A14").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B14").Value = "FPQ (Food propensity questionnair without portion sizes)"
$ws2.Range("C14").Value = 0

$ws2.Range("A15").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B15").Value = "FFQ (Food frequency questionnaire"
$ws2.Range("C15").Value = 1

$ws2.Range("A16").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B16").Value = "24HDR (24-h dietary recall"
$ws2.Range("C16").Value = 2

$ws2.Range("A17").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B17").Value = "3_d_FR_w (3-day weighing food record)"
$ws2.Range("C17").Value = 3

# Row 19 ("7_d_FR_w") is entered before row 18 ("7_d_FR") on purpose: this
# mirrors the original authoring order (visible in the shared-string table
# ordering) where the "...FR_w" label was typed before the "...FR" one even
# though it ended up one row below it.
$ws2.Range("A19").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B19").Value = "7_d_FR_w (7-day weighing food record)"
$ws2.Range("C19").Value = 5

$ws2.Range("A18").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B18").Value = "7_d_FR (7-day  food record; described portion sizes)"
$ws2.Range("C18").Value = 4

$ws2.Range("A20").Value = "DIETARY_ASSESS_INSTR"
$ws2.Range("B20").Value = "24HFL_FFQ (24-h short food list combined with FFQ"
$ws2.Range("C20").Value = 6

# --- Leave the selection/active sheet the way the saved file shows it:
#     Categories scrolled/selected over the new block, but Variables is the
#     sheet left active with the new row selected. ---
$ws2.Select()
$ws2.Range("A14:C20").Select()

$ws1.Select()
$ws1.Range("A129:D129").Select()
